# Updated cryptos list with GitHub Actions
# Note: a leading apostrophe forces Excel to store the value as literal text,
# preventing auto-conversion of numeric-looking strings (prices, percentages)
# into numbers, matching the original workbook where these columns are text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.606.19"
$ws.Range("E2").Value = "'  -1.60%  "

$ws.Range("D3").Value = "'1.666.51"
$ws.Range("E3").Value = "'  -3.29%  "

$ws.Range("E4").Value = "'  +0.01%  "

$ws.Range("D5").Value = "'215.06"
$ws.Range("E5").Value = "'  -1.83%  "

$ws.Range("E6").Value = "'  -2.42%  "

$ws.Range("E7").Value = "'  +0.03%  "

$ws.Range("D8").Value = "'23.87"
$ws.Range("E8").Value = "'  -2.04%  "

$ws.Range("D9").Value = "'0.262"
$ws.Range("E9").Value = "'  -0.67%  "

$ws.Range("E10").Value = "'  -1.71%  "

$ws.Range("D11").Value = "'0.0878"
$ws.Range("E11").Value = "'  -2.20%  "

$ws.Range("D12").Value = "'1.902.32"
$ws.Range("E12").Value = "'  -3.25%  "

$ws.Range("D13").Value = "'1.667.64"
$ws.Range("E13").Value = "'  -3.32%  "

$ws.Range("E14").Value = "'  -3.35%  "

$ws.Range("E15").Value = "'  -0.53%  "

$ws.Range("D16").Value = "'66.45"
$ws.Range("E16").Value = "'  -1.61%  "

$ws.Range("D17").Value = "'27.596.16"

$ws.Range("D18").Value = "'242.81"
$ws.Range("E18").Value = "'  +0.12%  "

$ws.Range("D19").Value = "'0.0" + [char]8323 + "0730"
$ws.Range("E19").Value = "'  -3.15%  "

$ws.Range("D20").Value = "'7.66"
$ws.Range("E20").Value = "'  -4.31%  "

$ws.Range("E21").Value = "'  -0.02%  "

$ws.Range("D22").Value = "'4.49"

$ws.Range("E23").Value = "'  -3.67%  "

$ws.Range("E24").Value = "'  -4.74%  "

$ws.Range("D25").Value = "'147.03"
$ws.Range("E25").Value = "'  -1.16%  "

$ws.Range("D26").Value = "'7.20"
$ws.Range("E26").Value = "'  -3.83%  "

$ws.Range("D27").Value = "'16.46"
$ws.Range("E27").Value = "'  -1.56%  "

$ws.Range("E29").Value = "'  -2.36%  "

$ws.Range("E30").Value = "'  +2.73%  "

$ws.Range("D31").Value = "'0.0502"
$ws.Range("E31").Value = "'  -1.51%  "

$ws.Range("E32").Value = "'  -2.47%  "

$ws.Range("D33").Value = "'1.471.16"
$ws.Range("E33").Value = "'  -1.44%  "

$ws.Range("E34").Value = "'  -4.71%  "

$ws.Range("E35").Value = "'  -5.16%  "

$ws.Range("E36").Value = "'  -1.25%  "

$ws.Range("D37").Value = "'0.929"
$ws.Range("E37").Value = "'  -2.69%  "

$ws.Range("D38").Value = "'0.576"
$ws.Range("E38").Value = "'  -5.07%  "

$ws.Range("E39").Value = "'  -1.55%  "

$ws.Range("D40").Value = "'69.44"
$ws.Range("E40").Value = "'  -1.67%  "

$ws.Range("E41").Value = "'  -4.34%  "

$ws.Range("E42").Value = "'  -0.01%  "

$ws.Range("E43").Value = "'  -3.00%  "

$ws.Range("E44").Value = "'  -7.60%  "

$ws.Range("D45").Value = "'1.810.07"
$ws.Range("E45").Value = "'  -3.17%  "

$ws.Range("D46").Value = "'0.787"
$ws.Range("E46").Value = "'  -1.75%  "

$ws.Range("E47").Value = "'  -2.47%  "

$ws.Range("D48").Value = "'89.23"
$ws.Range("E48").Value = "'  -1.98%  "

$ws.Range("E49").Value = "'  -4.06%  "

$ws.Range("E50").Value = "'  -2.05%  "

$ws.Range("D51").Value = "'7.89"
$ws.Range("E51").Value = "'  -4.59%  "
